$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.157.72"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "3.335.72"
$ws.Range("E3").Value = "  -5.02%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "604.14"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "163.91"
$ws.Range("E6").Value = "  -6.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -5.35%  "
$ws.Range("D9").Value = "3.327.46"
$ws.Range("E9").Value = "  -5.09%  "
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "6.70"
$ws.Range("E11").Value = "  -7.10%  "
$ws.Range("E12").Value = "  -8.77%  "
$ws.Range("D13").Value = "41.84"
$ws.Range("E13").Value = "  -9.60%  "
$ws.Range("E14").Value = "  -6.31%  "
$ws.Range("D15").Value = "3.875.72"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").Value = "7.73"
$ws.Range("E16").Value = "  -6.62%  "
$ws.Range("D17").Value = "68.121.86"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").Value = "3.339.93"
$ws.Range("E18").Value = "  -5.48%  "
$ws.Range("D19").Value = "559.99"
$ws.Range("E19").Value = "  -8.25%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "16.18"
$ws.Range("E21").Value = "  -6.62%  "
$ws.Range("E22").Value = "  -8.08%  "
$ws.Range("D23").Value = "8.46"
$ws.Range("E23").Value = "  -6.40%  "
$ws.Range("D24").Value = "89.98"
$ws.Range("E24").Value = "  -7.72%  "
$ws.Range("D25").Value = "14.28"
$ws.Range("E25").Value = "  -8.42%  "
$ws.Range("E26").Value = "  -6.40%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "31.22"
$ws.Range("E28").Value = "  -8.06%  "
$ws.Range("E29").Value = "  -10.61%  "
$ws.Range("D30").Value = "8.21"
$ws.Range("E30").Value = "  -8.51%  "
$ws.Range("E31").Value = "  -9.34%  "
$ws.Range("E32").Value = "  -6.95%  "
$ws.Range("D33").Value = "2.66"
$ws.Range("E33").Value = "  -10.96%  "
$ws.Range("D34").Value = "578.22"
$ws.Range("E34").Value = "  -9.29%  "
$ws.Range("E35").Value = "  -9.00%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "55.09"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").Value = "0.0451"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("E39").Value = "  -7.15%  "
$ws.Range("D40").Value = "0.0907"
$ws.Range("E40").Value = "  -8.82%  "
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").Value = "2.94"
$ws.Range("E42").Value = "  -18.27%  "
$ws.Range("D43").Value = "3.079.13"
$ws.Range("E43").Value = "  -8.22%  "
$ws.Range("D44").Value = "2.67"
$ws.Range("E44").Value = "  -7.98%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.278"
$ws.Range("E45").Value = "  -9.92%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0626"
$ws.Range("E46").Value = "  -15.86%  "
$ws.Range("D47").Value = "29.01"
$ws.Range("E47").Value = "  -9.78%  "
$ws.Range("D48").Value = "2.27"
$ws.Range("E48").Value = "  -11.16%  "
$ws.Range("D49").Value = "0.120"
$ws.Range("E49").Value = "  -7.09%  "
$ws.Range("E51").Value = "  -0.02%  "
